# Formulario.xlsx update:
#  - Se agrega "Tipo de visitante:" y se separa "Nombre del aprendiz:" en
#    "Nombre:" y "Apellidos:"
#  - Se agregan los campos "Número de celular" / "Correo electrónico" (ya
#    existian, ahora quedan mas abajo) y se agregan filas en blanco al final
#    del formulario para la funcion de calendario (fecha).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlContinuous = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$xlNone       = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$xlEdgeLeft   = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft
$xlEdgeTop    = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop
$xlEdgeBottom = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom
$xlEdgeRight  = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight

function Set-Label($row, $text, $sides) {
    $c = $ws.Cells.Item($row, 1)
    $c.Value = $text
    $c.Font.Name = "Arial"
    $c.Font.Size = 11
    $c.Font.Bold = $true

    if ($sides.Contains("L")) { $c.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous } else { $c.Borders.Item($xlEdgeLeft).LineStyle = $xlNone }
    if ($sides.Contains("R")) { $c.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous } else { $c.Borders.Item($xlEdgeRight).LineStyle = $xlNone }
    if ($sides.Contains("T")) { $c.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous } else { $c.Borders.Item($xlEdgeTop).LineStyle = $xlNone }
    if ($sides.Contains("B")) { $c.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous } else { $c.Borders.Item($xlEdgeBottom).LineStyle = $xlNone }
}

function Set-Input($row, $sides) {
    $c = $ws.Cells.Item($row, 2)
    $c.Font.Name = "Arial"
    $c.Font.Size = 11
    $c.Font.Bold = $false

    if ($sides.Contains("L")) { $c.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous } else { $c.Borders.Item($xlEdgeLeft).LineStyle = $xlNone }
    if ($sides.Contains("R")) { $c.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous } else { $c.Borders.Item($xlEdgeRight).LineStyle = $xlNone }
    if ($sides.Contains("T")) { $c.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous } else { $c.Borders.Item($xlEdgeTop).LineStyle = $xlNone }
    if ($sides.Contains("B")) { $c.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous } else { $c.Borders.Item($xlEdgeBottom).LineStyle = $xlNone }
}

# --- Columna A/B: etiquetas y campos del formulario ------------------------
Set-Label 1  "Tipo de visitante:"                  "TB"
Set-Input 1  "LRB"

Set-Label 2  "Nombre:"                              "LTB"
Set-Input 2  "LRTB"
$ws.Cells.Item(2, 2).Font.Size = 10

Set-Label 3  "Apellidos:"                           "LTB"
Set-Input 3  "LRTB"
$ws.Cells.Item(3, 2).Font.Size = 10

Set-Label 4  "Tipo de documento:"                   "TB"
Set-Input 4  "LRTB"

Set-Label 5  "Número de documento:"                 "TB"
Set-Input 5  "LRTB"

Set-Label 6  "Nombre del programa de formación:"    "TB"
Set-Input 6  "LRTB"

Set-Label 7  "Ficha:"                                ""
Set-Input 7  "LRTB"

Set-Label 8  "Centro:"                               "TB"
Set-Input 8  "LRTB"

Set-Label 9  "Número de celular"                     "TB"
Set-Input 9  "LRTB"

Set-Label 10 "Correo electrónico"                    "TB"
Set-Input 10 "LRTB"

# --- Filas en blanco adicionales (espacio para el calendario / firma) -----
for ($row = 11; $row -le 16; $row++) {
    $a = $ws.Cells.Item($row, 1)
    $b = $ws.Cells.Item($row, 2)
    $a.Borders.LineStyle = $xlNone
    $b.Borders.LineStyle = $xlNone
    $a.Font.Name = "Arial"
    $b.Font.Name = "Arial"
}

# Celda usada por el nuevo control de calendario (selector de fecha) en D13
$d13 = $ws.Cells.Item(13, 4)
$d13.Font.Name = "Arial"
$d13.Font.Size = 10
$d13.Borders.LineStyle = $xlNone

$ws.Range("D13").Select()
